# Commit: "Add bugs and DONE US."
#
# 1) Rename the header of column D on "Backlog" from "User story" to
#    "User story/Tâches technique".
# 2) Append a new "DONE" user story row (Maven / Ruahatu) to "Backlog".
# 3) Add a new "Bugs" worksheet after "Backlog" listing three bugs.

$wb = $excel.ActiveWorkbook
$backlog = $wb.Worksheets.Item("Backlog")

# --- 1) Header text tweak -------------------------------------------------
$backlog.Range("D1").Value = "User story/Tâches technique"

# --- 2) New row 52 on Backlog ---------------------------------------------
$backlog.Range("A52").Value = 51
$backlog.Range("B52").Value = "Maven"
$backlog.Range("C52").Value = "Ruahatu"
$backlog.Range("D52").Value = "Utiliser le module maven"
$backlog.Range("H52").Value = "DONE"

# Column B got wider (best-fit) to accommodate the longer theme labels.
$backlog.Columns.Item(2).ColumnWidth = 31.6

# Move the visible selection/scroll position like the author left it.
$backlog.Range("H53").Select()

# --- 3) New "Bugs" worksheet ------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$bugs = $wb.Worksheets.Add($null, $last)
$bugs.Name = "Bugs"

# Base look: bold/filled header row, plain data rows (same named styles the
# "Backlog" sheet already uses), with the date columns formatted afterwards.
$bugs.Range("A1:E1").Style = "Vérification"
$bugs.Range("A2:E4").Style = "Sortie"

$bugs.Range("B1:C1").NumberFormat = "dd/mm/yy;@"
$bugs.Range("B2:C4").NumberFormat = "dd/mm/yy;@"

$bugs.Range("A1").Value = "Numéro"
$bugs.Range("B1").Value = "Date saisie"
$bugs.Range("C1").Value = "Projet"
$bugs.Range("D1").Value = "Titre"
$bugs.Range("E1").Value = "Commentaire"

$bugs.Range("A2").Value = 1
$bugs.Range("B2").Value = 40559
$bugs.Range("C2").Value = "Ruahatu - Services"
$bugs.Range("D2").Value = "Services Web pas assez respectueux des principes REST"

$bugs.Range("A3").Value = 2
$bugs.Range("B3").Value = 40559
$bugs.Range("C3").Value = "Ruahatu - Client"
$bugs.Range("D3").Value = "Affichage - Liste poissons pas assez testé (Selenium)"

$bugs.Range("A4").Value = 3
$bugs.Range("B4").Value = 40559
$bugs.Range("C4").Value = "Ruahatu - Client"
$bugs.Range("D4").Value = "Authentification pas assez testé"

$bugs.Columns.Item(3).ColumnWidth = 17.28515625
$bugs.Columns.Item(4).ColumnWidth = 50.7109375
$bugs.Columns.Item(5).ColumnWidth = 53.140625

$bugs.Range("B4").Select()

# Leave "Backlog" as the selected/active tab, as in the authored workbook.
$backlog.Activate()
$backlog.Range("H53").Select()
